# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.675.43"
$ws.Range("E2").Value = "  -2.75%  "
$ws.Range("D3").Value = "3.555.77"
$ws.Range("E3").Value = "  -3.42%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "618.20"
$ws.Range("E5").Value = "  -6.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.52"
$ws.Range("E6").Value = "  -3.35%  "
$ws.Range("D7").Value = "3.550.65"
$ws.Range("E7").Value = "  -3.37%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.486"
$ws.Range("E9").Value = "  -2.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.141"
$ws.Range("E10").Value = "  -2.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.88"
$ws.Range("E11").Value = "  -3.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.432"
$ws.Range("E12").Value = "  -1.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000224"
$ws.Range("E13").Value = "  -3.26%  "
$ws.Range("D14").Value = "4.158.99"
$ws.Range("E14").Value = "  -3.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "32.25"
$ws.Range("E15").Value = "  -1.21%  "
$ws.Range("D16").Value = "3.560.66"
$ws.Range("E16").Value = "  -3.12%  "
$ws.Range("D17").Value = "67.726.50"
$ws.Range("E17").Value = "  -2.65%  "
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "15.63"
$ws.Range("E19").Value = "  -2.38%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.39"
$ws.Range("E20").Value = "  -0.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "454.50"
$ws.Range("E21").Value = "  -2.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.48"
$ws.Range("E22").Value = "  -2.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.642"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.64"
$ws.Range("E24").Value = "  -2.59%  "
$ws.Range("D25").Value = "3.703.38"
$ws.Range("E25").Value = "  -3.25%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  -6.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.59"
$ws.Range("E28").Value = "  -2.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.34"
$ws.Range("E29").Value = "  -6.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.56"
$ws.Range("E30").Value = "  -3.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.63"
$ws.Range("E31").Value = "  -2.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "25.97"
$ws.Range("E33").Value = "  -2.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.91"
$ws.Range("E34").Value = "  -4.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.159"
$ws.Range("E35").Value = "  -1.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.24"
$ws.Range("E36").Value = "  -3.14%  "
$ws.Range("D37").Value = "3.557.12"
$ws.Range("E37").Value = "  -3.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.07"
$ws.Range("E38").Value = "  -3.91%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "176.77"
$ws.Range("E41").Value = "  -1.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0884"
$ws.Range("E42").Value = "  -1.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.61"
$ws.Range("E43").Value = "  -6.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.08"
$ws.Range("E44").Value = "  -5.65%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.889"
$ws.Range("E45").Value = "  -4.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "29.00"
$ws.Range("E46").Value = "  +6.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.02"
$ws.Range("E47").Value = "  -1.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.59"
$ws.Range("E48").Value = "  -5.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.70"
$ws.Range("E49").Value = "  -1.34%  "
$ws.Range("E50").Value = "  -6.37%  "
$ws.Range("E51").Value = "  -3.91%  "
